$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.254.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.91%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.620.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.26%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'605.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.80%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'180.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.33%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '  +1.75%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.620.44"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.30%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '  +14.92%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '  +0.49%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '  +4.06%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '  +1.29%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("E15").Value = '  +5.02%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = '  +7.84%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'71.255.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.99%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.613.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.12%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'379.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.27%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'11.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.75%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'7.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.02%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '  +0.74%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'71.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.54%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = '  +6.23%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.00%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.86%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.756.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.43%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.996"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.25%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '  +7.09%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'531.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.71%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '  +3.23%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '  +5.59%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '  +3.68%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'165.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.90%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '  -1.13%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = '  +4.97%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = '  +8.45%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = '  +1.86%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '  +4.83%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'5.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.10%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'2.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.25%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '  +2.89%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'40.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.30%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D48").Value = "'3.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.57%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '  +5.24%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.531"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.62%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '  +7.37%  '
$ws.Range("E51").Style = "Normal"
